# The height of the "Implementierung des Frontends" task row (row 11) had
# drifted away from the standard 30pt row height used throughout the rest
# of the Gantt chart (it was sitting at an un-set/auto 30.75pt). Explicitly
# fix its height to the correct value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).RowHeight = 28.5

# Leave the selection where the user ended up after adjusting the row.
$ws.Range("L11").Select() | Out-Null
